$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-11, per diff
$newValues = @{
    2  = 3
    3  = 2
    4  = 5
    5  = 5
    6  = 3
    7  = 4
    8  = 2
    9  = 1
    10 = 2
    11 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
